# "what is the problem with L?"
#
# Adds a new "only binning" worksheet (JD / V[mag] / error_V[mag] columns,
# copied from the "after beaning" sheet's JD, V and error_V columns) right
# after "after beaning", makes it the active sheet/tab, and restores the
# various sheet selections that shifted around as a result.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "only binning" sheet right after "after beaning"
# ------------------------------------------------------------------
$afterBeaning = $wb.Worksheets.Item("after beaning")
$onlyBinning = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterBeaning)
$onlyBinning.Name = "only binning"

# Copy the bold/bordered header-cell formatting used elsewhere in the
# workbook (e.g. "after beaning"!C1) onto the new header row.
$afterBeaning.Range("C1").Copy()
$onlyBinning.Range("A1:C1").PasteSpecial(-4122)  # xlPasteFormats

$onlyBinning.Range("A1").Value = "JD - 2457651.0[day]"
$onlyBinning.Range("B1").Value = "V[mag]"
$onlyBinning.Range("C1").Value = "error_V[mag]"

# Data rows 2-16: JD, V[mag], error_V[mag] -- the same values as columns
# C, E, F of "after beaning" rows 2-16 (the un-binned / "after beaning" data).
$data = @(
    @(0.69693183999999997, 19.399999999999999, 0),
    @(0.7361548,            19.03,              0.08),
    @(0.73452949000000001,  19.11,              0.13),
    @(0.73629635999999998,  19.02,              0.13),
    @(0.73793518000000002,  18.98,              0.12),
    @(0.74201744999999997,  18.82,              0.07000000000000001),
    @(0.73999994999999996,  18.75,              0.12),
    @(0.74130099999999999,  18.75,              0.08),
    @(0.74262035000000004,  18.79,              0.1),
    @(0.74414860999999999,  18.79,              0.15),
    @(0.74817347000000001,  18.48,              0.05),
    @(0.74595498999999998,  18.46,              0.08),
    @(0.74763661999999997,  18.600000000000001, 0.09),
    @(0.74904883,           18.55,              0.09),
    @(0.75049745999999995,  18.34,              0.08)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $onlyBinning.Cells.Item($row, 1).Value = $data[$i][0]
    $onlyBinning.Cells.Item($row, 2).Value = $data[$i][1]
    $onlyBinning.Cells.Item($row, 3).Value = $data[$i][2]
}

# ------------------------------------------------------------------
# 2. Restore per-sheet selections that moved around in the source file
# ------------------------------------------------------------------

# "after beaning": columns C, E and F selected, landing on F1.
$afterBeaning.Activate()
$afterBeaning.Range("F:F").Select()

# "no beaning": B103 selected.
$noBeaning = $wb.Worksheets.Item("no beaning")
$noBeaning.Activate()
$noBeaning.Range("B103").Select()

# "pros data all": C40 selected.
$prosDataAll = $wb.Worksheets.Item("pros data all")
$prosDataAll.Activate()
$prosDataAll.Range("C40").Select()

# "the graph data": A44 selected.
$theGraphData = $wb.Worksheets.Item("the graph data")
$theGraphData.Activate()
$theGraphData.Range("A44").Select()

# ------------------------------------------------------------------
# 3. "only binning" ends up the active tab, cursor on A3.
# ------------------------------------------------------------------
$onlyBinning.Activate()
$onlyBinning.Range("A3").Select()
